# Januari-2021.xlsx — "10 years Finalization data"
#
# The original workbook has a single sheet ("Data Harian - Table") that
# contains a small header block (rows 1-5), a blank gap, and the daily
# weather data table at A9:K40 (a header row at row 9 followed by 31 daily
# rows, one per day of January).
#
# The commit adds a second worksheet named "Sheet1" that is a plain copy of
# that A9:K40 data table (re-based to A1:K32), reusing the same shared
# strings and the same header/data cell styles. The new sheet becomes the
# active / selected tab; the original sheet's own selection is left sitting
# on the data table it still contains.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data Harian - Table")

# --- add the new worksheet right after the existing one -------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$newSheet.Name = "Sheet1"

# --- copy the daily data table (header row + 31 day rows) over ------------
$src = $ws1.Range("A9:K40")
$src.Copy($newSheet.Range("A1"))

# Data rows on the new sheet wrap onto two lines (the sheet has no explicit
# column widths, unlike the source table), so give them the taller,
# auto-fit-sized row height; the header row keeps the single-line default.
for ($r = 2; $r -le 32; $r++) {
    $newSheet.Rows.Item($r).RowHeight = 28.8
}

# --- view state: keep the source sheet's selection on its table, and make
#     the new sheet the active / selected tab with its whole table selected
$ws1.Activate()
$ws1.Range("A9:K40").Select()

$newSheet.Activate()
$newSheet.Range("A1:K32").Select()
